# Auto-generated script to apply Valefor Profits market-data refresh
# (chore: update Sheets via scheduled runner)
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value2 = 3171.0908
$ws.Range("I64").Value2 = 3047.75
$ws.Range("K64").Value2 = 3047.75
$ws.Range("M64").Value2 = -2799.75
$ws.Range("H67").Value2 = 3171.0908
$ws.Range("I67").Value2 = 3047.75
$ws.Range("K67").Value2 = 3047.75
$ws.Range("M67").Value2 = -2189.75
$ws.Range("H76").Value2 = 3074.2856
$ws.Range("I76").Value2 = 2680
$ws.Range("J76").Value2 = 4750
$ws.Range("K76").Value2 = 2680
$ws.Range("L76").Value2 = 4750
$ws.Range("M76").Value2 = -2365
$ws.Range("N76").Value2 = -5380
$ws.Range("H79").Value2 = 3074.2856
$ws.Range("I79").Value2 = 2680
$ws.Range("J79").Value2 = 4750
$ws.Range("K79").Value2 = 2680
$ws.Range("L79").Value2 = 4750
$ws.Range("M79").Value2 = -1588
$ws.Range("N79").Value2 = -6934

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 2032.6136
$ws.Range("J61").Value2 = 2870
$ws.Range("L61").Value2 = 2870
$ws.Range("N61").Value2 = -3294
$ws.Range("H63").Value2 = 2402.0547
$ws.Range("I63").Value2 = 2358.0881
$ws.Range("K63").Value2 = 2358.0881
$ws.Range("M63").Value2 = -1672.0881
$ws.Range("H66").Value2 = 2402.0547
$ws.Range("I66").Value2 = 2358.0881
$ws.Range("K66").Value2 = 11790.4405
$ws.Range("M66").Value2 = -8358.440500000001
$ws.Range("H88").Value2 = 1756.2222
$ws.Range("I88").Value2 = 1000
$ws.Range("J88").Value2 = 1850.75
$ws.Range("K88").Value2 = 1000
$ws.Range("L88").Value2 = 1850.75
$ws.Range("M88").Value2 = -594
$ws.Range("N88").Value2 = -2662.75
$ws.Range("H91").Value2 = 1756.2222
$ws.Range("I91").Value2 = 1000
$ws.Range("J91").Value2 = 1850.75
$ws.Range("K91").Value2 = 1000
$ws.Range("L91").Value2 = 1850.75
$ws.Range("M91").Value2 = 404
$ws.Range("N91").Value2 = -4658.75
$ws.Range("H132").Value2 = 1637.1091
$ws.Range("I132").Value2 = 1160.58
$ws.Range("J132").Value2 = 6402.4
$ws.Range("K132").Value2 = 3481.74
$ws.Range("L132").Value2 = 19207.2
$ws.Range("M132").Value2 = -951.7399999999998
$ws.Range("N132").Value2 = -24267.2
$ws.Range("H136").Value2 = 2032.6136
$ws.Range("J136").Value2 = 2870
$ws.Range("L136").Value2 = 8610
$ws.Range("N136").Value2 = -13710

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 11907565
$ws.Range("I31").Value2 = 19231820
$ws.Range("J31").Value2 = 5651.75
$ws.Range("K31").Value2 = 19231820
$ws.Range("L31").Value2 = 5651.75
$ws.Range("M31").Value2 = -19231525
$ws.Range("N31").Value2 = -6241.75
$ws.Range("H34").Value2 = 11907565
$ws.Range("I34").Value2 = 19231820
$ws.Range("J34").Value2 = 5651.75
$ws.Range("K34").Value2 = 19231820
$ws.Range("L34").Value2 = 5651.75
$ws.Range("M34").Value2 = -19231618
$ws.Range("N34").Value2 = -6055.75
$ws.Range("H58").Value2 = 1470.4423
$ws.Range("I58").Value2 = 932.35
$ws.Range("J58").Value2 = 3264.0833
$ws.Range("K58").Value2 = 932.35
$ws.Range("L58").Value2 = 3264.0833
$ws.Range("M58").Value2 = -729.35
$ws.Range("N58").Value2 = -3670.0833
$ws.Range("H62").Value2 = 3211.7778
$ws.Range("I62").Value2 = 2850
$ws.Range("J62").Value2 = 3257
$ws.Range("K62").Value2 = 2850
$ws.Range("L62").Value2 = 3257
$ws.Range("M62").Value2 = -2226
$ws.Range("N62").Value2 = -4505
$ws.Range("H65").Value2 = 3211.7778
$ws.Range("I65").Value2 = 2850
$ws.Range("J65").Value2 = 3257
$ws.Range("K65").Value2 = 14250
$ws.Range("L65").Value2 = 16285
$ws.Range("M65").Value2 = -11130
$ws.Range("N65").Value2 = -22525
$ws.Range("H132").Value2 = 3279.0667
$ws.Range("I132").Value2 = 2698.3333
$ws.Range("J132").Value2 = 3666.2222
$ws.Range("K132").Value2 = 8094.999899999999
$ws.Range("L132").Value2 = 10998.6666
$ws.Range("M132").Value2 = -5564.999899999999
$ws.Range("N132").Value2 = -16058.6666
$ws.Range("H136").Value2 = 1470.4423
$ws.Range("I136").Value2 = 932.35
$ws.Range("J136").Value2 = 3264.0833
$ws.Range("K136").Value2 = 2797.05
$ws.Range("L136").Value2 = 9792.249899999999
$ws.Range("M136").Value2 = -247.0500000000002
$ws.Range("N136").Value2 = -14892.2499

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value2 = 1151.7826
$ws.Range("I11").Value2 = 165.77777
$ws.Range("J11").Value2 = 1785.6428
$ws.Range("K11").Value2 = 497.33331
$ws.Range("L11").Value2 = 5356.928400000001
$ws.Range("M11").Value2 = -357.33331
$ws.Range("N11").Value2 = -5636.928400000001

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value2 = 27000
$ws.Range("J46").Value2 = 27000
$ws.Range("L46").Value2 = 27000
$ws.Range("N46").Value2 = -27312
$ws.Range("H55").Value2 = 1733.3334
$ws.Range("I55").Value2 = 1733.3334
$ws.Range("J55").Value2 = 0
$ws.Range("K55").Value2 = 1733.3334
$ws.Range("L55").Value2 = 0
$ws.Range("M55").Value2 = -1406.3334
$ws.Range("N55").ClearContents()
$ws.Range("H70").Value2 = 4731.8
$ws.Range("I70").Value2 = 4412.5
$ws.Range("K70").Value2 = 4412.5
$ws.Range("M70").Value2 = -4142.5
$ws.Range("H73").Value2 = 4731.8
$ws.Range("I73").Value2 = 4412.5
$ws.Range("K73").Value2 = 4412.5
$ws.Range("M73").Value2 = -3476.5
$ws.Range("H80").Value2 = 56121.367
$ws.Range("I80").Value2 = 2455.5557
$ws.Range("J80").Value2 = 104420.6
$ws.Range("K80").Value2 = 2455.5557
$ws.Range("L80").Value2 = 104420.6
$ws.Range("M80").Value2 = -1457.5557
$ws.Range("N80").Value2 = -106416.6
$ws.Range("H83").Value2 = 56121.367
$ws.Range("I83").Value2 = 2455.5557
$ws.Range("J83").Value2 = 104420.6
$ws.Range("K83").Value2 = 12277.7785
$ws.Range("L83").Value2 = 522103
$ws.Range("M83").Value2 = -7285.7785
$ws.Range("N83").Value2 = -532087
$ws.Range("H132").Value2 = 11501854
$ws.Range("I132").Value2 = 19615260
$ws.Range("J132").Value2 = 7860.4165
$ws.Range("K132").Value2 = 58845780
$ws.Range("L132").Value2 = 23581.2495
$ws.Range("M132").Value2 = -58843250
$ws.Range("N132").Value2 = -28641.2495

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value2 = 2685.1516
$ws.Range("I132").Value2 = 1872.48
$ws.Range("J132").Value2 = 5224.75
$ws.Range("K132").Value2 = 5617.440000000001
$ws.Range("L132").Value2 = 15674.25
$ws.Range("M132").Value2 = -3087.440000000001
$ws.Range("N132").Value2 = -20734.25

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 1419.5424
$ws.Range("I132").Value2 = 1155.3489
$ws.Range("J132").Value2 = 2129.5625
$ws.Range("K132").Value2 = 3466.0467
$ws.Range("L132").Value2 = 6388.6875
$ws.Range("M132").Value2 = -936.0466999999999
$ws.Range("N132").Value2 = -11448.6875
